# Using ordinal attribute values to calculate Moran's I:
# split the previous single "whole sample" total row into a
# "first half" subtotal row (new row 31, rows 2:15) and keep the
# "all rows" grand total (shifts down to row 32), which now also
# naturally folds in the new subtotal row 16.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Autocorrelation")
$ws.Select()

# Insert a new row before row 31, shifting rows 31-35 down to 32-36
$ws.Rows.Item(31).Insert()

# Add SUM formulas into row 16 for columns C, D, E
$ws.Range("C16").Formula = "=SUM(C2:C15)"
$ws.Range("D16:E16").Formula = "=SUM(D2:D15)"

# Add SUM formulas into the newly inserted row 31 for columns C, D, E
$ws.Range("C31").Formula = "=SUM(C17:C30)"
$ws.Range("D31:E31").Formula = "=SUM(D17:D30)"

# Select D33 to match final cursor position
$ws.Range("D33").Select()
